# Add a new "numpy" worksheet after "Spaltenvektor" (the last sheet) and
# populate it with the ">20 filter" example, mirroring the original
# author's commit ("Added xlsx to gitignore" accompanied a larger content
# change that introduced this sheet).

$wb = $excel.ActiveWorkbook

# --- create the new sheet at the end of the workbook -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "numpy"

# Template cells already present in the workbook that carry the exact
# formatting (fonts/fills/borders) we need to reuse for the new sheet -
# copying them keeps the style table free of duplicate / near-duplicate
# cellXfs entries.
$txn = $wb.Worksheets.Item("transactions")
$loop1 = $wb.Worksheets.Item("forloop1")

# --- header row (bold, centered, bottom border) -------------------------
# ("Tag" is entered first so the shared-string table ends up in the same
#  order the original authoring session produced it.)
$txn.Range("C5").Copy($ws.Range("B5"))
$ws.Range("B5").Value = "Tag"

# --- title row ---------------------------------------------------------
$ws.Range("B2").Value = "Wichtige Analyse"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Size = 24

$txn.Range("C5").Copy($ws.Range("C5"))
$ws.Range("C5").Value = "Werte"

$txn.Range("C5").Copy($ws.Range("D5"))
$ws.Range("D5").Value = ">20"

# --- data rows (6-17): Tag index, raw value, IF(value>20, value, 0) -----
$values = @(22, 31, 34, 16, 8, 23, 7, 26, 45, 31, 3, 35)
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 6 + $i
    $loop1.Range("B4").Copy($ws.Range("B" + $row))
    $ws.Range("B" + $row).Value = $i + 1

    $loop1.Range("B4").Copy($ws.Range("C" + $row))
    $ws.Range("C" + $row).Value = $values[$i]

    $loop1.Range("B4").Copy($ws.Range("D" + $row))
    $ws.Range("D" + $row).Formula = "=IF(C" + $row + ">20,C" + $row + ",0)"
}

# stray note cell next to the second data row
$ws.Range("G7").Value = "dfdf"

# --- totals row ----------------------------------------------------------
$txn.Range("E12").Copy($ws.Range("B19"))
$ws.Range("B19").Value = "Gesamt"

$txn.Range("E12").Copy($ws.Range("C19"))
$ws.Range("C19").Formula = "=+SUM(C6:C17)"

$txn.Range("E12").Copy($ws.Range("D19"))
$ws.Range("D19").Formula = "=+SUM(D6:D17)"

# --- view state: make "numpy" the active/visible tab, restore the
#     previously-active sheet's selection to where the user last left it --
$spalten = $wb.Worksheets.Item("Spaltenvektor")
$spalten.Activate()
$spalten.Range("H24").Select()

$ws.Activate()
$excel.ActiveWindow.Zoom = 88
$ws.Range("K29").Select()
